# Apply the "New Orleans shard 189" edit:
#
#   - The sheet that was "hotel_info" (1st tab) becomes "review_info" and
#     is reduced to just its header row (24 review columns + the leading
#     "STR" column).
#   - The sheet that was "review_info" (2nd tab, header-only) becomes
#     "hotel_info" and receives the hotel_info header + data row, with a
#     new "State" column inserted right after "Hotel_Name" (value
#     "Louisiana" for this row).

$wb = $excel.ActiveWorkbook

$sheetA = $wb.Worksheets.Item(1)   # currently "hotel_info"   -> becomes "review_info"
$sheetB = $wb.Worksheets.Item(2)   # currently "review_info"  -> becomes "hotel_info"

# sheetB currently holds the (header-only) review_info table; wipe it
# before dropping the hotel_info data in, so no stale cells are left
# dangling past column J.
$sheetB.Cells.Clear()

# --- Move the hotel_info header/data (A:B stay put, C:I shift right one
#     column to make room for the new "State" column at C) onto sheetB,
#     using Range.Copy so the original cell types (numbers vs. text, e.g.
#     the text-typed "503"/"17"/"515") are preserved verbatim. ---
$sheetA.Range("A1:B2").Copy($sheetB.Range("A1"))
$sheetA.Range("C1:I2").Copy($sheetB.Range("D1"))

# New "State" column.
$sheetB.Cells.Item(1, 3).Value = "State"
$sheetB.Cells.Item(2, 3).Value = "Louisiana"

# --- Replace sheetA's old hotel_info content with the review_info header
#     row (no data rows). ---
$sheetA.Cells.Clear()

$reviewHeaders = @(
    "STR","reviewer_ID","reviewer_name","Review_ID","Date_of_scraping","ReviewURL",
    "Tripadvisor_gcode","Tripadvisor_dcode","Tripadvisor_rcode","review_date","review_title",
    "review_content","review_rating","trip_month","trip_purpose","value","rooms","Location",
    "Cleanliness","Sleep Quality","Service","Picture(yes=1)","respondent","response_date","response_text"
)
for ($i = 0; $i -lt $reviewHeaders.Length; $i++) {
    $sheetA.Cells.Item(1, $i + 1).Value = $reviewHeaders[$i]
}

# --- Rename through a temporary name so the two tabs never collide while
#     swapping ("hotel_info" -> "review_info" and vice versa at once). ---
$sheetA.Name = "__tmp_shard189__"
$sheetB.Name = "hotel_info"
$sheetA.Name = "review_info"
